$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.098888
$ws.Range("H2").Value = 6.296664
$ws.Range("I2").Value = 0.1082453658858517
$ws.Range("J2").Value = 0.1082453658858517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 4.879420661690667
$ws.Range("R2").Value = 43.914785955216
$ws.Range("S2").Value = 0.004842815401780955
$ws.Range("T2").Value = 0.004842815401780953
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.098888
$ws.Range("H3").Value = 6.296664
$ws.Range("I3").Value = 0.1082453658858517
$ws.Range("J3").Value = 0.1082453658858517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 40.45589758933067
$ws.Range("R3").Value = 364.103078303976
$ws.Range("S3").Value = 0.04015239872157266
$ws.Range("T3").Value = 0.04015239872157265
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.098888
$ws.Range("H4").Value = 6.296664
$ws.Range("I4").Value = 0.1082453658858517
$ws.Range("J4").Value = 0.1082453658858517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 63.72823900153334
$ws.Range("R4").Value = 573.5541510138
$ws.Range("S4").Value = 0.06325015176249807
$ws.Range("T4").Value = 0.06325015176249806
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.87514366666667
$ws.Range("H5").Value = 47.625431
$ws.Range("I5").Value = 0.8187243600843848
$ws.Range("J5").Value = 0.8187243600843847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.324764666666667
$ws.Range("N5").Value = 6.974294
$ws.Range("O5").Value = 0.04473923998638302
$ws.Range("P5").Value = 0.04473923998638301
$ws.Range("Q5").Value = 36.90597307452378
$ws.Range("R5").Value = 332.153757670714
$ws.Range("S5").Value = 0.03662910562851315
$ws.Range("T5").Value = 0.03662910562851315
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.87514366666667
$ws.Range("H6").Value = 47.625431
$ws.Range("I6").Value = 0.8187243600843848
$ws.Range("J6").Value = 0.8187243600843847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.27491966666667
$ws.Range("N6").Value = 57.824759
$ws.Range("O6").Value = 0.3709387315842666
$ws.Range("P6").Value = 0.3709387315842665
$ws.Range("Q6").Value = 305.9921188717921
$ws.Range("R6").Value = 2753.929069846129
$ws.Range("S6").Value = 0.303696575646842
$ws.Range("T6").Value = 0.3036965756468419
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.87514366666667
$ws.Range("H7").Value = 47.625431
$ws.Range("I7").Value = 0.8187243600843848
$ws.Range("J7").Value = 0.8187243600843847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.5843220284293504
$ws.Range("P7").Value = 0.5843220284293504
$ws.Range("Q7").Value = 482.0147381723139
$ws.Range("R7").Value = 4338.132643550825
$ws.Range("S7").Value = 0.4783986788090296
$ws.Range("T7").Value = 0.4783986788090295
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.416064
$ws.Range("H8").Value = 4.248192
$ws.Range("I8").Value = 0.07303027402976368
$ws.Range("J8").Value = 0.07303027402976367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.324764666666667
$ws.Range("N8").Value = 6.974294
$ws.Range("O8").Value = 0.04473923998638302
$ws.Range("P8").Value = 0.04473923998638301
$ws.Range("Q8").Value = 3.292015552938667
$ws.Range("R8").Value = 29.628139976448
$ws.Range("S8").Value = 0.003267318956088912
$ws.Range("T8").Value = 0.003267318956088912
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.416064
$ws.Range("H9").Value = 4.248192
$ws.Range("I9").Value = 0.07303027402976368
$ws.Range("J9").Value = 0.07303027402976367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.27491966666667
$ws.Range("N9").Value = 57.824759
$ws.Range("O9").Value = 0.3709387315842666
$ws.Range("P9").Value = 0.3709387315842665
$ws.Range("Q9").Value = 27.29451984285866
$ws.Range("R9").Value = 245.650678585728
$ws.Range("S9").Value = 0.02708975721585195
$ws.Range("T9").Value = 0.02708975721585194
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.416064
$ws.Range("H10").Value = 4.248192
$ws.Range("I10").Value = 0.07303027402976368
$ws.Range("J10").Value = 0.07303027402976367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.5843220284293504
$ws.Range("P10").Value = 0.5843220284293504
$ws.Range("Q10").Value = 42.99575062293333
$ws.Range("R10").Value = 386.9617556064
$ws.Range("S10").Value = 0.04267319785782283
$ws.Range("T10").Value = 0.04267319785782282
